$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 87, pushing existing rows 87-130 down to 89-132.
$ws.Rows.Item(87).Insert()
$ws.Rows.Item(87).Insert()

# Populate new row 87: Packham's Triumph / Primera
$ws.Cells.Item(87,1).Value  = 4
$ws.Cells.Item(87,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(87,3).Value  = "Los Lagos"
$ws.Cells.Item(87,4).Value  = 44466
$ws.Cells.Item(87,5).Value  = 10
$ws.Cells.Item(87,6).Value  = "Fruta"
$ws.Cells.Item(87,7).Value  = 100104
$ws.Cells.Item(87,8).Value  = "Frutos de pepita"
$ws.Cells.Item(87,9).Value  = 100104005
$ws.Cells.Item(87,10).Value = "Pera"
$ws.Cells.Item(87,11).Value = "Packham's Triumph"
$ws.Cells.Item(87,12).Value = "Primera"
$ws.Cells.Item(87,13).Value = 200
$ws.Cells.Item(87,14).Value = 16000
$ws.Cells.Item(87,15).Value = 16000
$ws.Cells.Item(87,16).Value = 16000
$ws.Cells.Item(87,17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(87,18).Value = "Región de O'Higgins"
$ws.Cells.Item(87,19).Value = 1067
$ws.Cells.Item(87,20).Value = 15

# Populate new row 88: Packham's Triumph / Segunda
$ws.Cells.Item(88,1).Value  = 4
$ws.Cells.Item(88,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(88,3).Value  = "Los Lagos"
$ws.Cells.Item(88,4).Value  = 44466
$ws.Cells.Item(88,5).Value  = 10
$ws.Cells.Item(88,6).Value  = "Fruta"
$ws.Cells.Item(88,7).Value  = 100104
$ws.Cells.Item(88,8).Value  = "Frutos de pepita"
$ws.Cells.Item(88,9).Value  = 100104005
$ws.Cells.Item(88,10).Value = "Pera"
$ws.Cells.Item(88,11).Value = "Packham's Triumph"
$ws.Cells.Item(88,12).Value = "Segunda"
$ws.Cells.Item(88,13).Value = 100
$ws.Cells.Item(88,14).Value = 13000
$ws.Cells.Item(88,15).Value = 13000
$ws.Cells.Item(88,16).Value = 13000
$ws.Cells.Item(88,17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(88,18).Value = "Región de O'Higgins"
$ws.Cells.Item(88,19).Value = 867
$ws.Cells.Item(88,20).Value = 15
